# Version 9 : Adding Data Set (Data Drivern)
#
# The "Test Steps" sheet gets a new "Data_Set" column inserted right
# before the existing "Results" column (which shifts from F to G).
# The username/password step-level keywords collapse from two distinct
# action keywords (input_uname / input_upass) into one generic "input"
# keyword, with the actual test data now living in the new Data_Set column.

$wb = $excel.ActiveWorkbook

$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsCases = $wb.Worksheets.Item("Test Cases")

# --- Test Steps sheet -------------------------------------------------

# Insert a new column F (pushes the old F "Results"/"PASS" column to G).
$wsSteps.Columns.Item(6).Insert()

# Give the new column its header width/name.
$wsSteps.Columns.Item(6).ColumnWidth = 38.57

$wsSteps.Range("F1").Value = "Data_Set"

# Row 2 / 9 : openBrowser step -> chrome data
$wsSteps.Range("F2").Value = "chrome"
$wsSteps.Range("F9").Value = "chrome"

# Row 4 / 11 : Enter username -> generic "input" keyword + email data set
$wsSteps.Range("E4").Value = "input"
$wsSteps.Range("F4").Value = "rajasingh.nadar@infosys.com.vmstdemo"
$wsSteps.Range("E11").Value = "input"
$wsSteps.Range("F11").Value = "rajasingh.nadar@infosys.com.vmstdemo"

# Row 5 / 12 : Enter password -> generic "input" keyword + password data set
$wsSteps.Range("E5").Value = "input"
$wsSteps.Range("F5").Value = "Raja@1506`$`$`$`$`$`$"
$wsSteps.Range("E12").Value = "input"
$wsSteps.Range("F12").Value = "Raja@1506`$`$`$`$`$`$"

# Update the view's selection + dimension to match the post-edit state.
$wsSteps.Range("E17").Select()

# --- Test Cases sheet ---------------------------------------------------
# No content changes; only the remembered selection moved.
# (Selecting this range last keeps "Test Cases" as the active/visible tab,
# matching the workbook's activeTab bookmark.)
$wsCases.Range("E12").Select()
